# Updates cryptos.xlsx symbol data (price/volume/hour columns) to the
# latest scrape snapshot. Values are written with a leading apostrophe so
# Excel stores them as text (matching the original inlineStr cells) rather
# than auto-converting numeric-looking strings (e.g. "331.36", "-0.39%",
# "18") into real numbers/percentages. Style is reset to "Normal" right
# after so the quote-prefix formatting Excel applies does not linger on
# the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'331.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.39%"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "'18"
$ws.Range("G2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'41.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'1.31%"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Value = "'18"
$ws.Range("G3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'5.663"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.97%"
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").Value = "'18"
$ws.Range("G4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'0.08351"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.45%"
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").Value = "'18"
$ws.Range("G5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'8.788"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.42%"
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").Value = "'18"
$ws.Range("G6").Style = "Normal"

# Row 7
$ws.Range("D7").Value = "'1.994"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-3.08%"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").Value = "'18"
$ws.Range("G7").Style = "Normal"

# Row 8
$ws.Range("D8").Value = "'4.467"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.67%"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "'18"
$ws.Range("G8").Style = "Normal"

# Row 9
$ws.Range("D9").Value = "'2.910"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-2.98%"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").Value = "'18"
$ws.Range("G9").Style = "Normal"

# Row 10
$ws.Range("D10").Value = "'0.9254"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.22%"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Value = "'18"
$ws.Range("G10").Style = "Normal"

# Row 11
$ws.Range("D11").Value = "'0.1291"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.02%"
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").Value = "'18"
$ws.Range("G11").Style = "Normal"

# Row 12
$ws.Range("D12").Value = "'0.1968"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.67%"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").Value = "'18"
$ws.Range("G12").Style = "Normal"

# Row 13
$ws.Range("D13").Value = "'0.09432"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'2.61%"
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").Value = "'18"
$ws.Range("G13").Style = "Normal"

# Row 14
$ws.Range("D14").Value = "'0.03881"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'4.75%"
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").Value = "'18"
$ws.Range("G14").Style = "Normal"

# Row 15
$ws.Range("E15").Value = "'0.95%"
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").Value = "'18"
$ws.Range("G15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'0.001297"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.02%"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").Value = "'18"
$ws.Range("G16").Style = "Normal"

# Row 17
$ws.Range("D17").Value = "'0.006106"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.84%"
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").Value = "'18"
$ws.Range("G17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'3.442"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.97%"
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").Value = "'18"
$ws.Range("G18").Style = "Normal"

# Row 19
$ws.Range("D19").Value = "'0.3535"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.54%"
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").Value = "'18"
$ws.Range("G19").Style = "Normal"

# Row 20
$ws.Range("D20").Value = "'8.465"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-2.48%"
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").Value = "'18"
$ws.Range("G20").Style = "Normal"

# Row 21
$ws.Range("D21").Value = "'0.1372"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-3.47%"
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").Value = "'18"
$ws.Range("G21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'0.2479"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-6.77%"
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").Value = "'18"
$ws.Range("G22").Style = "Normal"

# Row 23
$ws.Range("D23").Value = "'0.04405"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.74%"
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").Value = "'18"
$ws.Range("G23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'0.001274"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.01%"
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").Value = "'18"
$ws.Range("G24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'0.004379"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.80%"
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").Value = "'18"
$ws.Range("G25").Style = "Normal"

# Row 26
$ws.Range("E26").Value = "'-1.82%"
$ws.Range("E26").Style = "Normal"
$ws.Range("G26").Value = "'18"
$ws.Range("G26").Style = "Normal"

# Row 27
$ws.Range("G27").Value = "'18"
$ws.Range("G27").Style = "Normal"

# Row 28
$ws.Range("G28").Value = "'18"
$ws.Range("G28").Style = "Normal"

# Row 29
$ws.Range("G29").Value = "'18"
$ws.Range("G29").Style = "Normal"

# Row 30
$ws.Range("G30").Value = "'18"
$ws.Range("G30").Style = "Normal"

# Row 31
$ws.Range("G31").Value = "'18"
$ws.Range("G31").Style = "Normal"

# Row 32
$ws.Range("G32").Value = "'18"
$ws.Range("G32").Style = "Normal"

# Row 33
$ws.Range("G33").Value = "'18"
$ws.Range("G33").Style = "Normal"

# Row 34
$ws.Range("G34").Value = "'18"
$ws.Range("G34").Style = "Normal"

# Row 35
$ws.Range("G35").Value = "'18"
$ws.Range("G35").Style = "Normal"

# Row 36
$ws.Range("G36").Value = "'18"
$ws.Range("G36").Style = "Normal"

# Row 37
$ws.Range("G37").Value = "'18"
$ws.Range("G37").Style = "Normal"

# Row 38
$ws.Range("G38").Value = "'18"
$ws.Range("G38").Style = "Normal"

# Row 39
$ws.Range("D39").Value = "'0.02841"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-0.99%"
$ws.Range("E39").Style = "Normal"
$ws.Range("G39").Value = "'18"
$ws.Range("G39").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.05539"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.95%"
$ws.Range("E40").Style = "Normal"
$ws.Range("G40").Value = "'18"
$ws.Range("G40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.007934"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'2.06%"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").Value = "'18"
$ws.Range("G41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'0.1436"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.37%"
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").Value = "'18"
$ws.Range("G42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'0.009307"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-6.22%"
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").Value = "'18"
$ws.Range("G43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.002240"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'1.13%"
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").Value = "'18"
$ws.Range("G44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'0.01111"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.81%"
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").Value = "'18"
$ws.Range("G45").Style = "Normal"

# Row 46
$ws.Range("D46").Value = "'0.00007270"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'6.67%"
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").Value = "'18"
$ws.Range("G46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.22%"
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").Value = "'18"
$ws.Range("G47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'0.003445"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'14.87%"
$ws.Range("E48").Style = "Normal"
$ws.Range("G48").Value = "'18"
$ws.Range("G48").Style = "Normal"

# Row 49
$ws.Range("D49").Value = "'0.002279"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.21%"
$ws.Range("E49").Style = "Normal"
$ws.Range("G49").Value = "'18"
$ws.Range("G49").Style = "Normal"

# Row 50
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.22%"
$ws.Range("E50").Style = "Normal"
$ws.Range("G50").Value = "'18"
$ws.Range("G50").Style = "Normal"

# Row 51
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.22%"
$ws.Range("E51").Style = "Normal"
$ws.Range("G51").Value = "'18"
$ws.Range("G51").Style = "Normal"
